# Applies the "Automatic update of files." commit:
#   - every data row's "Förändrad" (column C) timestamp bumps from 46064 to 46065
#   - the case rows (rows 6-95, excluding a handful that already sit in the
#     right spot) get re-shuffled into a new order; each case keeps its full
#     row of data (columns A-Z) intact as it moves to its new row position.
#
# The mapping below gives, for each destination row (2..95, in order), which
# source row (in the *current*, pre-edit sheet) supplies that row's data.
# Rows that don't move simply map to themselves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 95

# destination row (index 0 => row 2, index 1 => row 3, ...) -> source row number
$srcForDest = @(2,3,4,5,7,6,8,9,10,11,12,13,14,15,17,16,18,20,19,24,21,22,23,25,26,27,28,29,30,31,32,38,35,40,41,33,36,39,34,42,73,37,76,92,59,95,43,94,80,70,93,58,87,91,69,60,61,50,83,64,89,66,46,57,85,67,62,51,86,90,71,56,47,81,75,84,44,72,49,68,45,65,88,74,82,48,63,54,52,53,55,78,77,79)

$rowCount = $lastRow - $firstRow + 1

# Snapshot the whole table (A:Z) before we start overwriting anything.
$srcRange = $ws.Range("A$firstRow`:Z$lastRow")
$origValues   = $srcRange.Value()
$origFormulas = $srcRange.Formula()

$lastCol = 26  # Z

# Build the destination arrays by pulling each destination row's data from
# its mapped source row (1-based within the snapshot arrays).
$newValues   = New-Object 'object[,]' $rowCount,$lastCol
$newFormulas = New-Object 'object[,]' $rowCount,$lastCol

for ($i = 0; $i -lt $rowCount; $i++) {
    $srcRowNum = $srcForDest[$i]
    $srcIdx = $srcRowNum - $firstRow + 1
    for ($c = 1; $c -le $lastCol; $c++) {
        $newValues[$i, $c-1]   = $origValues[$srcIdx, $c]
        $newFormulas[$i, $c-1] = $origFormulas[$srcIdx, $c]
    }
}

# Write back plain values for the non-formula columns (A..R = cols 1..18).
$valRange = $ws.Range("A$firstRow`:R$lastRow")
$valOut = New-Object 'object[,]' $rowCount,18
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($c = 1; $c -le 18; $c++) {
        $valOut[$i, $c-1] = $newValues[$i, $c-1]
    }
}
$valRange.Value = $valOut

# Write back formulas for the hyperlink columns (S..Z = cols 19..26).
$formRange = $ws.Range("S$firstRow`:Z$lastRow")
$formOut = New-Object 'object[,]' $rowCount,8
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($c = 19; $c -le 26; $c++) {
        $f = $newFormulas[$i, $c-1]
        if ($null -eq $f -or $f -eq "") {
            $formOut[$i, $c-19] = ""
        } else {
            $formOut[$i, $c-19] = $f
        }
    }
}
$formRange.Formula = $formOut

# Finally, bump the "Förändrad" column (C) to the new date serial for every
# data row (this applies uniformly, independent of the reshuffle above).
$ws.Range("C$firstRow`:C$lastRow").Value = 46065
